{"js": "// Release QR update: bump UVVM library version numbers and the\n// \"Last update\" footer date (see commit message: \"Release: updated\n// CHANGES.TXT and QRs with release version numbers\").\n//\n//   - \"UVVM Utility Library (UVVM-Util), version 2.10.0 and up\"\n//         -> \"... version 2.11.0 and up\"\n//   - \"UVVM VVC Framework, version 2.7.0 and up\"\n//         -> \"... version 2.7.1 and up\"\n//   - footer \"... Last update: 2019-12-03\"\n//         -> \"... Last update: 2020-01-23\"\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// --- 1) UVVM Utility Library version: 2.10.0 -> 2.11.0 -------------------\nconst utilPara = paragraphs.items.find(\n  (p) => p.text.indexOf(\"UVVM Utility Library\") !== -1 && p.text.indexOf(\"version\") !== -1\n);\nif (!utilPara) {\n  throw new Error(\"Could not find the 'UVVM Utility Library' version paragraph\");\n}\nconst utilVersionHits = utilPara.search(\"10\", { matchCase: true });\nutilVersionHits.load(\"items\");\nawait context.sync();\nif (utilVersionHits.items.length === 0) {\n  throw new Error(\"Could not find the '10' version token to update\");\n}\nutilVersionHits.items[0].insertText(\"11\", \"Replace\");\nawait context.sync();\n\n// --- 2) UVVM VVC Framework version: 2.7.0 -> 2.7.1 ------------------------\nconst vvcPara = paragraphs.items.find((p) => p.text.indexOf(\"UVVM VVC Framework, version\") !== -1);\nif (!vvcPara) {\n  throw new Error(\"Could not find the 'UVVM VVC Framework' version paragraph\");\n}\n// Insert \".1\" right after \"2.7\" (mirrors how Word records a mid-text\n// insertion); the trailing \".0 and up\" run (after the _GoBack bookmark)\n// then loses its leading \".0\" in a second, separate step.\nconst majorMinorHits = vvcPara.search(\"2.7\", { matchCase: true });\nmajorMinorHits.load(\"items\");\nawait context.sync();\nif (majorMinorHits.items.length === 0) {\n  throw new Error(\"Could not find the '2.7' version token to update\");\n}\nconst afterSeven = majorMinorHits.items[0].getRange(\"End\");\nafterSeven.insertText(\".1\", \"Before\");\nawait context.sync();\n\nconst trailingHits = vvcPara.search(\".0 and up\", { matchCase: true });\ntrailingHits.load(\"items\");\nawait context.sync();\nif (trailingHits.items.length === 0) {\n  throw new Error(\"Could not find the trailing '.0 and up' token to update\");\n}\ntrailingHits.items[0].insertText(\" and up\", \"Replace\");\nawait context.sync();\n\n// --- 3) Footer \"Last update\" date: 2019-12-03 -> 2020-01-23 ---------------\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nfor (const section of sections.items) {\n  const footer = section.getFooter(\"Primary\");\n  const dateHits = footer.search(\"2019-12-03\", { matchCase: true });\n  dateHits.load(\"items\");\n  await context.sync();\n  for (const hit of dateHits.items) {\n    hit.insertText(\"2020-01-23\", \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Release QR update: bump UVVM library version numbers and the\n# \"Last update\" footer date (see commit message: \"Release: updated\n# CHANGES.TXT and QRs with release version numbers\").\n#\n#   - \"UVVM Utility Library (UVVM-Util), version 2.10.0 and up\"\n#         -> \"... version 2.11.0 and up\"\n#   - \"UVVM VVC Framework, version 2.7.0 and up\"\n#         -> \"... version 2.7.1 and up\"\n#   - footer \"... Last update: 2019-12-03\"\n#         -> \"... Last update: 2020-01-23\"\n\n$d = $word.ActiveDocument\n\n# --- 1) UVVM Utility Library version: 2.10.0 -> 2.11.0 ---------------------\n$utilFound = $false\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text\n    if (($text -like \"*UVVM Utility Library*\") -and ($text -like \"*version*\")) {\n        $r = $p.Range\n        if ($r.Find.Execute(\"10\")) {\n            # Replace just the trailing \"0\" of \"10\" with \"1\" so \"10\" -> \"11\",\n            # leaving the surrounding runs untouched.\n            $zero = $d.Range($r.End - 1, $r.End)\n            $zero.Text = \"1\"\n            $utilFound = $true\n        }\n        break\n    }\n}\nif (-not $utilFound) {\n    throw \"Could not find the 'UVVM Utility Library' version token to update\"\n}\n\n# --- 2) UVVM VVC Framework version: 2.7.0 -> 2.7.1 --------------------------\n$vvcFound = $false\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text\n    if ($text -like \"*UVVM VVC Framework, version*\") {\n        $r = $p.Range\n        if ($r.Find.Execute(\"2.7\")) {\n            # Insert \".1\" right after \"2.7\" (keeps the _GoBack bookmark in place).\n            $insertPoint = $d.Range($r.End, $r.End)\n            $insertPoint.InsertAfter(\".1\")\n        } else {\n            throw \"Could not find the '2.7' version token to update\"\n        }\n\n        $r2 = $p.Range\n        if ($r2.Find.Execute(\".0 and up\")) {\n            $r2.Text = \" and up\"\n            $vvcFound = $true\n        } else {\n            throw \"Could not find the trailing '.0 and up' token to update\"\n        }\n        break\n    }\n}\nif (-not $vvcFound) {\n    throw \"Could not find the 'UVVM VVC Framework' version paragraph\"\n}\n\n# --- 3) Footer \"Last update\" date: 2019-12-03 -> 2020-01-23 -----------------\n$dateFound = $false\nfor ($i = 1; $i -le $d.Sections.Count; $i++) {\n    $section = $d.Sections($i)\n    $footer = $section.Footers(1)  # wdHeaderFooterPrimary\n    $fr = $footer.Range\n    while ($fr.Find.Execute(\"2019-12-03\")) {\n        $fr.Text = \"2020-01-23\"\n        $dateFound = $true\n        $fr = $footer.Range\n        $fr.Start = $fr.End\n    }\n}\nif (-not $dateFound) {\n    throw \"Could not find the '2019-12-03' footer date to update\"\n}\n"}
